$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-20 22:21:12"
$wsZh.Range("H4").Value = "2016-03-20 22:21:36"
$wsZh.Range("E5").Value = "2016-03-20 22:21:12"
$wsZh.Range("H5").Value = "2016-03-20 22:21:36"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-20 22:21:15"
$wsDe.Range("H4").Value = "2016-03-20 22:21:42"
$wsDe.Range("E5").Value = "2016-03-20 22:21:15"
$wsDe.Range("H5").Value = "2016-03-20 22:21:42"
